$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H64").Value = 9076
$ws.Range("I64").Value = 8247
$ws.Range("J64").Value = 10181.333
$ws.Range("K64").Value = 8247
$ws.Range("L64").Value = 10181.333
$ws.Range("M64").Value = -7999
$ws.Range("N64").Value = -10677.333
$ws.Range("H67").Value = 9076
$ws.Range("I67").Value = 8247
$ws.Range("J67").Value = 10181.333
$ws.Range("K67").Value = 8247
$ws.Range("L67").Value = 10181.333
$ws.Range("M67").Value = -7389
$ws.Range("N67").Value = -11897.333
$ws.Range("H106").Value = 6027.4443
$ws.Range("I106").Value = 6027.4443
$ws.Range("K106").Value = 6027.4443
$ws.Range("M106").Value = -5396.4443
$ws.Range("H107").Value = 352
$ws.Range("I107").Value = 341.77777
$ws.Range("K107").Value = 341.77777
$ws.Range("M107").Value = 1578.22223
$ws.Range("H138").Value = 3853.28
$ws.Range("J138").Value = 5350.7417
$ws.Range("L138").Value = 16052.2251
$ws.Range("N138").Value = -26332.2251

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 1950.4445
$ws.Range("I74").Value = 1950.4445
$ws.Range("K74").Value = 1950.4445
$ws.Range("M74").Value = -1076.4445
$ws.Range("H77").Value = 1950.4445
$ws.Range("I77").Value = 1950.4445
$ws.Range("K77").Value = 9752.2225
$ws.Range("M77").Value = -5384.2225
$ws.Range("H88").Value = 4018.818
$ws.Range("J88").Value = 4356.3335
$ws.Range("L88").Value = 4356.3335
$ws.Range("N88").Value = -5168.3335
$ws.Range("H91").Value = 4018.818
$ws.Range("J91").Value = 4356.3335
$ws.Range("L91").Value = 4356.3335
$ws.Range("N91").Value = -7164.3335
$ws.Range("H110").Value = 2066.7144
$ws.Range("I110").Value = 1344.8334
$ws.Range("K110").Value = 1344.8334
$ws.Range("M110").Value = 700.1666
$ws.Range("H132").Value = 3789.3125
$ws.Range("I132").Value = 3537.12
$ws.Range("J132").Value = 4690
$ws.Range("K132").Value = 10611.36
$ws.Range("L132").Value = 14070
$ws.Range("M132").Value = -8081.360000000001
$ws.Range("N132").Value = -19130

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 956.0714
$ws.Range("I107").Value = 952.6923
$ws.Range("K107").Value = 952.6923
$ws.Range("M107").Value = 967.3077
$ws.Range("H132").Value = 81119.375
$ws.Range("J132").Value = 81119.375
$ws.Range("L132").Value = 81119.375
$ws.Range("N132").Value = -91239.375
$ws.Range("H134").Value = 3641.3948
$ws.Range("I134").Value = 3572.9375
$ws.Range("J134").Value = 4006.5
$ws.Range("K134").Value = 10718.8125
$ws.Range("L134").Value = 12019.5
$ws.Range("M134").Value = -8183.8125
$ws.Range("N134").Value = -17089.5
$ws.Range("H140").Value = 117186.664
$ws.Range("J140").Value = 117186.664
$ws.Range("L140").Value = 117186.664
$ws.Range("N140").Value = -127546.664

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 59899.8
$ws.Range("J16").Value = 71999.5
$ws.Range("L16").Value = 71999.5
$ws.Range("N16").Value = -72573.5
$ws.Range("H62").Value = 4780.2666
$ws.Range("I62").Value = 4589.625
$ws.Range("J62").Value = 4998.143
$ws.Range("K62").Value = 4589.625
$ws.Range("L62").Value = 4998.143
$ws.Range("M62").Value = -3965.625
$ws.Range("N62").Value = -6246.143
$ws.Range("H65").Value = 4780.2666
$ws.Range("I65").Value = 4589.625
$ws.Range("J65").Value = 4998.143
$ws.Range("K65").Value = 22948.125
$ws.Range("L65").Value = 24990.715
$ws.Range("M65").Value = -19828.125
$ws.Range("N65").Value = -31230.715
$ws.Range("H86").Value = 15166158
$ws.Range("I86").Value = 18534618
$ws.Range("J86").Value = 8087
$ws.Range("K86").Value = 18534618
$ws.Range("L86").Value = 8087
$ws.Range("M86").Value = -18533495
$ws.Range("N86").Value = -10333
$ws.Range("H89").Value = 15166158
$ws.Range("I89").Value = 18534618
$ws.Range("J89").Value = 8087
$ws.Range("K89").Value = 92673090
$ws.Range("L89").Value = 40435
$ws.Range("M89").Value = -92667474
$ws.Range("N89").Value = -51667
$ws.Range("H105").Value = 4547.6665
$ws.Range("I105").Value = 1721.75
$ws.Range("J105").Value = 10199.5
$ws.Range("K105").Value = 1721.75
$ws.Range("L105").Value = 10199.5
$ws.Range("M105").Value = 25.25
$ws.Range("N105").Value = -13693.5
$ws.Range("H107").Value = 1097.65
$ws.Range("I107").Value = 1190.4
$ws.Range("K107").Value = 1190.4
$ws.Range("M107").Value = 729.5999999999999
$ws.Range("H113").Value = 59899.8
$ws.Range("J113").Value = 71999.5
$ws.Range("L113").Value = 71999.5
$ws.Range("N113").Value = -76339.5
$ws.Range("H122").Value = 103074.1
$ws.Range("I122").Value = 180349.47
$ws.Range("K122").Value = 541048.41
$ws.Range("M122").Value = -538598.41
$ws.Range("H134").Value = 1699.8334
$ws.Range("I134").Value = 1399.6666
$ws.Range("K134").Value = 4198.9998
$ws.Range("M134").Value = -1663.9998

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("L104").ClearContents()
$ws.Range("H60").Value = 430
$ws.Range("I60").Value = 578.3333
$ws.Range("K60").Value = 1734.9999
$ws.Range("M60").Value = -1483.9999
$ws.Range("H104").Value = 7029
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 7029
$ws.Range("K104").Value = 0
$ws.Range("M104").Value = 21087
$ws.Range("H140").Value = 2164.3157
$ws.Range("I140").Value = 1254.3529
$ws.Range("K140").Value = 3763.0587
$ws.Range("M140").Value = 1416.9413
$ws.Range("N104").Value = -26329

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 5590.6
$ws.Range("I80").Value = 3414.5715
$ws.Range("J80").Value = 7494.625
$ws.Range("K80").Value = 3414.5715
$ws.Range("L80").Value = 7494.625
$ws.Range("M80").Value = -2416.5715
$ws.Range("N80").Value = -9490.625
$ws.Range("H83").Value = 5590.6
$ws.Range("I83").Value = 3414.5715
$ws.Range("J83").Value = 7494.625
$ws.Range("K83").Value = 17072.8575
$ws.Range("L83").Value = 37473.125
$ws.Range("M83").Value = -12080.8575
$ws.Range("N83").Value = -47457.125
$ws.Range("H126").Value = 8765.333000000001
$ws.Range("I126").Value = 6898.25
$ws.Range("K126").Value = 20694.75
$ws.Range("M126").Value = -18224.75
$ws.Range("H132").Value = 3808
$ws.Range("I132").Value = 3851.0186
$ws.Range("J132").Value = 3614.4167
$ws.Range("K132").Value = 11553.0558
$ws.Range("L132").Value = 10843.2501
$ws.Range("M132").Value = -9023.0558
$ws.Range("N132").Value = -15903.2501

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("M23").ClearContents()
$ws.Range("H23").Value = 1342666.6
$ws.Range("I23").Value = 1342666.6
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1342666.6
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = -1342436.6
$ws.Range("H46").Value = 2906.125
$ws.Range("J46").Value = 2550
$ws.Range("L46").Value = 2550
$ws.Range("N46").Value = -2926
$ws.Range("H122").Value = 7332.4287
$ws.Range("I122").Value = 7512.615
$ws.Range("K122").Value = 22537.845
$ws.Range("M122").Value = -20087.845
$ws.Range("H132").Value = 4166.3335
$ws.Range("I132").Value = 4166.3335
$ws.Range("K132").Value = 12499.0005
$ws.Range("M132").Value = -9969.000499999998
$ws.Range("H136").Value = 4492.933
$ws.Range("I136").Value = 3470.75
$ws.Range("K136").Value = 10412.25
$ws.Range("M136").Value = -7862.25

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H29").Value = 13001.429
$ws.Range("I29").Value = 13127.5
$ws.Range("J29").Value = 12833.333
$ws.Range("K29").Value = 13127.5
$ws.Range("L29").Value = 12833.333
$ws.Range("M29").Value = -12837.5
$ws.Range("N29").Value = -13413.333
$ws.Range("H47").Value = 47943.89
$ws.Range("I47").Value = 43500
$ws.Range("K47").Value = 43500
$ws.Range("H122").Value = 2838.0667
$ws.Range("I122").Value = 3799.5
$ws.Range("K122").Value = 11398.5
$ws.Range("M122").Value = -8948.5
$ws.Range("H132").Value = 6563.522
$ws.Range("I132").Value = 4763.4116
$ws.Range("J132").Value = 11663.833
$ws.Range("K132").Value = 14290.2348
$ws.Range("L132").Value = 34991.499
$ws.Range("M132").Value = -11760.2348
$ws.Range("N132").Value = -40051.499
$ws.Range("H136").Value = 4083.5806
$ws.Range("I136").Value = 4027.64
$ws.Range("K136").Value = 12082.92
$ws.Range("M136").Value = -9532.92
$ws.Range("M47").Value = -42928
